$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 120, shifting existing rows
# 120-122 down to 121-123 (preserving their original values/formatting).
$ws.Rows.Item(120).Insert()

# Populate the newly inserted row 120 with this week's entry. Columns
# A, B, C, E, F, G, H, I, J, K, R are identical across every row in this
# sheet (same market/product classification), so copy them from the row
# directly below (the shifted former row120, now row121).
$ws.Range("A120").Value = $ws.Range("A121").Value2
$ws.Range("B120").Value = $ws.Range("B121").Value2
$ws.Range("C120").Value = $ws.Range("C121").Value2
$ws.Range("D120").Value = 45265
$ws.Range("E120").Value = $ws.Range("E121").Value2
$ws.Range("F120").Value = $ws.Range("F121").Value2
$ws.Range("G120").Value = $ws.Range("G121").Value2
$ws.Range("H120").Value = $ws.Range("H121").Value2
$ws.Range("I120").Value = $ws.Range("I121").Value2
$ws.Range("J120").Value = $ws.Range("J121").Value2
$ws.Range("K120").Value = $ws.Range("K121").Value2
$ws.Range("L120").Value = "Primera"
$ws.Range("M120").Value = 100
$ws.Range("N120").Value = 40000
$ws.Range("O120").Value = 40000
$ws.Range("P120").Value = 40000
$ws.Range("Q120").Value = "$/caja 16 kilos"
$ws.Range("R120").Value = $ws.Range("R121").Value2
$ws.Range("S120").Value = 2500
$ws.Range("T120").Value = 16
